# Update the "人气/浏览量" (column F) figures across all four sheets to the
# freshly generated gh-pages snapshot values.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { cell address -> new value }
$updates = @{
    "展览"     = @{
        "F3"  = 279
        "F4"  = 1163
        "F6"  = 2777
        "F9"  = 111
        "F12" = 709
        "F13" = 111
        "F15" = 1771
        "F18" = 207
        "F19" = 260
    }
    "演出"     = @{
        "F4"  = 29
        "F7"  = 25
        "F10" = 50
        "F12" = 55
        "F23" = 31
    }
    "本地生活" = @{
        "F3" = 799
        "F4" = 2033
        "F5" = 272
    }
    "全部类型" = @{
        "F3"  = 799
        "F4"  = 2033
        "F5"  = 272
        "F9"  = 29
        "F11" = 279
        "F12" = 1163
        "F15" = 25
        "F17" = 2777
        "F20" = 50
        "F22" = 55
        "F24" = 111
        "F28" = 709
        "F29" = 111
        "F32" = 1771
        "F37" = 207
        "F44" = 31
        "F45" = 260
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellAddr in $cellUpdates.Keys) {
        $ws.Range($cellAddr).Value = $cellUpdates[$cellAddr]
    }
}
